{"js": "// 1) Mark the first row of the confidence-rating table as a repeating\n//    header row (w:trPr/w:tblHeader on the first <w:tr>).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.headerRowCount = 1;\nawait context.sync();\n\n// 2) Trim the \"Sketch the t-distribution using the t-distribution applet.\"\n//    sentence off the end of the P-value bullet point, leaving just:\n//    \"Find the P-value and compare it to the level of significance.\"\nconst searchText =\n  \"Find the P-value and compare it to the level of significance. \" +\n  \"Sketch the t-distribution using the t-distribution applet.\";\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"Find the P-value and compare it to the level of significance.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# 1) Trim the \"Sketch the t-distribution using the t-distribution applet.\"\n#    sentence off the end of the P-value bullet point, leaving just:\n#    \"Find the P-value and compare it to the level of significance.\"\n$d = $word.ActiveDocument\n$target = \"Find the P-value and compare it to the level of significance. Sketch the t-distribution using the t-distribution applet.\"\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $t = $r.Text\n    if ($t -ne $null -and $t.TrimEnd([char]13, [char]7) -eq $target) {\n        $r.Text = \"Find the P-value and compare it to the level of significance.\"\n        break\n    }\n}\n\n# 2) Mark the first row of the confidence-rating table as a repeating\n#    header row (w:trPr/w:tblHeader on the first <w:tr>).\n$tbl = $d.Tables.Item(1)\n$tbl.Rows.Item(1).HeadingFormat = $true\n"}
